$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K5").Value = 18.025
$ws.Range("L5").Value = 26.32

$ws.Range("K6").Value = 24
$ws.Range("L6").Value = 63.38

$ws.Range("K7").Value = 107.733
$ws.Range("L7").Value = 81.62

$ws.Range("K8").Value = 164.758
$ws.Range("L8").Value = 198.67
